# 5.2.1.1a — add a new "Q" column (year 2023) mirroring the existing
# per-year columns (D..P = 2010..2022).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 is a thin "separator" row: every used column in it just carries
# the bottom-border style (s=4) with no value. Copy that style from P3.
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)   # xlPasteFormats

# Row 4 header: year value 2023, same style as the other year header
# cells (s=13).
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Q4").Value = 2023

# Row 5 is a blank spacer row under the section header; also gets taller
# (24.75 -> 27) to fit the new column's content.
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows("5:5").RowHeight = 27

# Data rows 6, 8, 9, 10: plain numeric values, same style as column P.
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("Q6").Value = 1209

$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)
$ws.Range("Q8").Value = 373

$ws.Range("P9").Copy()
$ws.Range("Q9").PasteSpecial(-4122)
$ws.Range("Q9").Value = 115

$ws.Range("P10").Copy()
$ws.Range("Q10").PasteSpecial(-4122)
$ws.Range("Q10").Value = 781

# Row 7: P7 already holds the "-" placeholder text; Q7 mirrors it
# (same style + same text).
$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$ws.Range("Q7").Value = "-"

$excel.CutCopyMode = $false
